$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $s = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $s
}

# Row 2
Set-TextValue $ws.Range('D2') '30.451.45'
Set-TextValue $ws.Range('E2') '  +0.22%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.868.05'
Set-TextValue $ws.Range('E3') '  -0.34%  '

# Row 4
Set-TextValue $ws.Range('D4') '1.000'
Set-TextValue $ws.Range('E4') '  +0.02%  '

# Row 5
Set-TextValue $ws.Range('D5') '235.41'
Set-TextValue $ws.Range('E5') '  -1.29%  '

# Row 6
Set-TextValue $ws.Range('E6') '  +0.00%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.4826'
Set-TextValue $ws.Range('E7') '  -0.06%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.2798'
Set-TextValue $ws.Range('E8') '  -1.00%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.06509'
Set-TextValue $ws.Range('E9') '  -0.37%  '

# Row 10
Set-TextValue $ws.Range('D10') '1.873.98'
Set-TextValue $ws.Range('E10') '  +0.32%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.07438'
Set-TextValue $ws.Range('E11') '  -0.49%  '

# Row 12
Set-TextValue $ws.Range('D12') '16.25'
Set-TextValue $ws.Range('E12') '  -1.71%  '

# Row 13
Set-TextValue $ws.Range('D13') '5.074'
Set-TextValue $ws.Range('E13') '  -0.50%  '

# Row 14
Set-TextValue $ws.Range('D14') '87.26'
Set-TextValue $ws.Range('E14') '  -1.07%  '

# Row 15
Set-TextValue $ws.Range('D15') '0.6417'
Set-TextValue $ws.Range('E15') '  -2.61%  '

# Row 16
Set-TextValue $ws.Range('D16') '30.438.61'
Set-TextValue $ws.Range('E16') '  +0.27%  '

# Row 17
Set-TextValue $ws.Range('D17') '0.9999'
Set-TextValue $ws.Range('E17') '  -0.02%  '

# Row 18
Set-TextValue $ws.Range('D18') '12.97'
Set-TextValue $ws.Range('E18') '  -2.75%  '

# Row 19
Set-TextValue $ws.Range('D19') '231.72'
Set-TextValue $ws.Range('E19') '  +3.91%  '

# Row 20
Set-TextValue $ws.Range('D20') '0.000007496'
Set-TextValue $ws.Range('E20') '  -1.80%  '

# Row 21
Set-TextValue $ws.Range('D21') '2.108.79'
Set-TextValue $ws.Range('E21') '  -0.24%  '

# Row 22
Set-TextValue $ws.Range('E22') '  +0.03%  '

# Row 23
Set-TextValue $ws.Range('D23') '5.147'
Set-TextValue $ws.Range('E23') '  -3.08%  '

# Row 24
Set-TextValue $ws.Range('B24') 'Chainlink'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D24') '6.100'
Set-TextValue $ws.Range('E24') '  -1.60%  '

# Row 25
Set-TextValue $ws.Range('B25') 'Cosmos'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D25') '9.308'
Set-TextValue $ws.Range('E25') '  +0.31%  '

# Row 26
Set-TextValue $ws.Range('B26') 'Monero'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D26') '167.85'
Set-TextValue $ws.Range('E26') '  +1.17%  '

# Row 27
Set-TextValue $ws.Range('B27') 'EthereumClassic'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D27') '18.41'
Set-TextValue $ws.Range('E27') '  -1.94%  '

# Row 28
Set-TextValue $ws.Range('B28') 'LidoDAOToken'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D28') '1.914'
Set-TextValue $ws.Range('E28') '  -3.38%  '

# Row 29
Set-TextValue $ws.Range('B29') 'Stellar'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D29') '0.1024'
Set-TextValue $ws.Range('E29') '  +8.75%  '

# Row 30
Set-TextValue $ws.Range('B30') 'Toncoin'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D30') '1.379'
Set-TextValue $ws.Range('E30') '  -5.45%  '

# Row 31
Set-TextValue $ws.Range('B31') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D31') '4.262'
Set-TextValue $ws.Range('E31') '  -1.51%  '

# Row 32
Set-TextValue $ws.Range('B32') 'Filecoin'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D32') '3.999'
Set-TextValue $ws.Range('E32') '  -0.60%  '

# Row 33
Set-TextValue $ws.Range('B33') 'Hedera'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D33') '0.04987'
Set-TextValue $ws.Range('E33') '  -1.43%  '

# Row 34
Set-TextValue $ws.Range('B34') 'ARBITRUM'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D34') '1.175'
Set-TextValue $ws.Range('E34') '  -3.52%  '

# Row 35
Set-TextValue $ws.Range('B35') 'ImmutableX'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D35') '0.7392'
Set-TextValue $ws.Range('E35') '  -1.87%  '

# Row 36
Set-TextValue $ws.Range('B36') 'Frax'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D36') '0.9997'
Set-TextValue $ws.Range('E36') '  +0.12%  '

# Row 37
Set-TextValue $ws.Range('B37') 'HuobiToken'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D37') '2.711'
Set-TextValue $ws.Range('E37') '  +0.41%  '

# Row 38
Set-TextValue $ws.Range('B38') 'VeChain'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D38') '0.01946'
Set-TextValue $ws.Range('E38') '  +5.63%  '

# Row 39
Set-TextValue $ws.Range('B39') 'MXToken'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D39') '2.636'
Set-TextValue $ws.Range('E39') '  +0.63%  '

# Row 40
Set-TextValue $ws.Range('B40') 'TrustWalletToken'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D40') '0.9188'
Set-TextValue $ws.Range('E40') '  +1.26%  '

# Row 41
Set-TextValue $ws.Range('B41') 'RenderToken'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D41') '2.050'
Set-TextValue $ws.Range('E41') '  -1.98%  '

# Row 42
Set-TextValue $ws.Range('B42') 'Quant'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D42') '105.95'
Set-TextValue $ws.Range('E42') '  -0.98%  '

# Row 43
Set-TextValue $ws.Range('B43') 'PaxDollar'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D43') '0.9963'
Set-TextValue $ws.Range('E43') '  -0.68%  '

# Row 44
Set-TextValue $ws.Range('B44') 'TheSandbox'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D44') '0.4197'
Set-TextValue $ws.Range('E44') '  -2.62%  '

# Row 45
Set-TextValue $ws.Range('B45') 'FraxShare'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D45') '5.571'
Set-TextValue $ws.Range('E45') '  -6.35%  '

# Row 46
Set-TextValue $ws.Range('B46') 'Aptos'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D46') '7.208'
Set-TextValue $ws.Range('E46') '  -3.42%  '

# Row 47
Set-TextValue $ws.Range('B47') 'Aave'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D47') '61.55'
Set-TextValue $ws.Range('E47') '  -5.10%  '

# Row 48
Set-TextValue $ws.Range('B48') 'Algorand'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D48') '0.1226'
Set-TextValue $ws.Range('E48') '  -5.86%  '

# Row 49
Set-TextValue $ws.Range('B49') 'EnergySwap'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '8.837'
Set-TextValue $ws.Range('E49') '  -2.63%  '

# Row 50
Set-TextValue $ws.Range('B50') 'NEARProtocol'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D50') '1.430'
Set-TextValue $ws.Range('E50') '  -4.48%  '

# Row 51
Set-TextValue $ws.Range('B51') 'Elrond'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range('D51') '33.48'
Set-TextValue $ws.Range('E51') '  -2.15%  '
